$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Format as text first so Excel doesn't auto-convert the dd-mm-yyyy
# looking strings into date serials (column A in this sheet stores
# these values as plain text, like all the preceding rows).
$ws.Range("A30:A31").NumberFormat = "@"

$ws.Range("A30").Value = "05-10-2021"
$ws.Range("B30").Value = 10000
$ws.Range("D30").Value = 0

$ws.Range("A31").Value = "06-10-2021"
$ws.Range("B31").Value = 10000
$ws.Range("D31").Value = 0

# Drop the temporary text format so the new cells end up with the
# same default (no explicit style) formatting as the existing rows.
$ws.Range("A30:A31").ClearFormats()
